$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(12).Delete()

$ws.Range("L1").Value = "PORTO_ID"
$ws.Range("M1").Value = "PORTO_NAMA"
$ws.Range("N1").Value = "INV_ID"
$ws.Range("O1").Value = "KODE"
$ws.Range("P1").Value = "NO_URUT"
$ws.Range("Q1").Value = "PREFIX_SURAT"
$ws.Range("R1").Value = "PORTO_ID_OJK"
$ws.Range("S1").Value = "KETERANGAN"
$ws.Range("T1").Value = "STATUS_REGISTER"
$ws.Range("U1").Value = "KETERANGAN_REGISTER"
